$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell E8 text value: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active-cell selection on the sheet (matches the added <selection> in sheetView)
$ws.Range("E8").Select()
